$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update the Eintragsdatum value (B10) from 17.12.2023 to 01.01.2025
$ws.Range("B10").Value = "01.01.2025"

# Reflect the new selection state captured in the saved file
$ws.Range("B14").Select()
